# Pico Dram Tester Timings.xlsx - apply commit edits
#
# 1) Notes sheet: insert two blank rows (splitting the numbered note list so
#    that note #2 and note #3 each sit in their own row with a blank
#    separator row above/between them), and fix the typo / wording in note #2
#    ("prgressively" -> "progressively", and clarify "nop instructions" ->
#    "relevant nop instructions").
# 2) On each of the six chip-timing sheets, fix the start-time formulas for
#    the two 'write' branch landing points (G13 / G15) which had been wired
#    to the wrong predecessor cells, and bump the 'Delays' cycle counts
#    (column B) that compensate for the corrected branch timing so that the
#    overall cycle counts used elsewhere stay consistent.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Notes sheet
# ---------------------------------------------------------------------
$notes = $wb.Worksheets.Item("Notes")

# Insert a blank row above old row 3 ("2.0" / ...) and another blank row
# above what becomes the "3.0" row, matching the target layout:
#   row2 = 1.0 ...
#   row3 = (blank)
#   row4 = 2.0 ...
#   row5 = (blank)
#   row6 = 3.0 ...
#   row7 = (trailing note)
$notes.Rows.Item(3).Insert()
$notes.Rows.Item(5).Insert()

# Fix the wording of note #2 (now in row 4, column B)
$notes.Range("B4").Value = "To use a sheet, the idea is to zero out the all the 'Delays' values and then progressively increase each one in order until all the timing cells between the relevant nop instructions are white. "

# ---------------------------------------------------------------------
# 2) Chip timing sheets
# ---------------------------------------------------------------------
$sheetNames = @("4164 - 100ns", "4164 - 120ns", "4164 - 150ns", "4816 - 100ns", "4816 - 120ns", "4816 - 150ns")

# New 'Delays' values (column B) per sheet. Only the rows that actually
# change value are listed; B8 changes on every sheet, B7 only changes on
# the 4816 sheets.
$delayUpdates = @{
    "4164 - 100ns" = @{ 8 = 12.0 }
    "4164 - 120ns" = @{ 8 = 17.0 }
    "4164 - 150ns" = @{ 8 = 24.0 }
    "4816 - 100ns" = @{ 7 = 5.0; 8 = 9.0 }
    "4816 - 120ns" = @{ 7 = 5.0; 8 = 14.0 }
    "4816 - 150ns" = @{ 7 = 5.0; 8 = 23.0 }
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Apply the corrected 'Delays' counts for this sheet.
    $rows = $delayUpdates[$name]
    foreach ($r in $rows.Keys) {
        $ws.Cells.Item($r, 2).Value = $rows[$r]
    }

    # G13 ('skip_wr:' label row) is reached directly from the row 10
    # 'jmp !x skip_wr' instruction, so its start time must come from H10,
    # not from the fall-through chain.
    $ws.Range("G13").Formula = "=H10"

    # G15 ('skip_wr2:' label row) is the convergence point of both the
    # write and non-write branches; its start time must come from H12
    # (the end of the write branch), not H14.
    $ws.Range("G15").Formula = "=H12"
}

$wb.Application.Calculate()
